$wb = $excel.ActiveWorkbook

# --- Metadata sheet: Count 7 -> 10 ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B21").Value = "'10"
# Re-apply the original cell style (text-forcing quote-prefix otherwise
# swaps in a new style); B20 shares the same style as B21.
$meta.Range("B20").Copy()
$meta.Range("B21").PasteSpecial(-4122)

# --- Concepts sheet: update existing rows & append new concept rows ---
$concepts = $wb.Worksheets.Item("Concepts")

# RL display text update
$concepts.Range("C3").Value = "In a relationship, not married, living with partner (Domestic Partner)"

# SD display text update
$concepts.Range("C6").Value = "Legally Separated"

# S display text update
$concepts.Range("C7").Value = "Never Married"

# New rows: Annuled, Interlocutory, Polygamous (Level is always "1")
$concepts.Range("A9").Value = "'1"
$concepts.Range("B9").Value = "A"
$concepts.Range("C9").Value = "Annuled"

$concepts.Range("A10").Value = "'1"
$concepts.Range("B10").Value = "I"
$concepts.Range("C10").Value = "Interlocutory"

$concepts.Range("A11").Value = "'1"
$concepts.Range("B11").Value = "P"
$concepts.Range("C11").Value = "Polygamous"

# Apply the same style as the rest of the data rows (style from row 8)
$concepts.Range("A8:D8").Copy()
$concepts.Range("A9:D11").PasteSpecial(-4122)
